$d = $word.ActiveDocument

# 1. Title heading and bold repeat of title (both occurrences replaced identically)
$d.Content.Find.Execute("Play Jolly Roger 2 Slot for Free - Enjoy Pirate-Themed Gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Play Jolly Roger 2 Slot Free - Exciting Pirate-Themed Game", 2)
$d.Content.Find.Execute("Play Jolly Roger 2 Slot for Free - Enjoy Pirate-Themed Gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Play Jolly Roger 2 Slot Free - Exciting Pirate-Themed Game", 2)

# 2. "What we like" bullet list - reorder/replace
$d.Content.Find.Execute("Exciting Bonus Quest Feature with four special levels", $true, $false, $false, $false, $false, $true, 1, $false, "Pirate-themed slot with exciting gameplay", 2)
$d.Content.Find.Execute("Excellent graphics and sound design for authentic pirate-themed game", $true, $false, $false, $false, $false, $true, 1, $false, "Randomly activated Compass feature with multipliers", 2)
$d.Content.Find.Execute("Medium volatility with a potential win of 5,000 times the bet", $true, $false, $false, $false, $false, $true, 1, $false, "Bonus Quest Feature with four special levels", 2)
$d.Content.Find.Execute("Compass feature adds a multiplier ranging from 2x to 10x", $true, $false, $false, $false, $false, $true, 1, $false, "Maximum potential win of 5,000 times the bet", 2)

# 3. "What we don't like" bullet list
$d.Content.Find.Execute("Fixed pay lines may limit betting options for some players", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting options for high rollers", 2)
$d.Content.Find.Execute("No progressive jackpot feature", $true, $false, $false, $false, $false, $true, 1, $false, "Higher volatility may not suit all players", 2)

# 4. Italic description paragraph
$d.Content.Find.Execute("Read our review of Jolly Roger 2, play for free, and enjoy exciting Bonus Quest Feature with four special levels and a potential win of 5,000 times the bet.", $true, $false, $false, $false, $false, $true, 1, $false, "Play Jolly Roger 2 for free and enjoy the excitement of a pirate-themed slot game with potential big payouts.", 2)
